{"js": "// Update the date label and all multiplication problems in the table.\n// Each source text is unique in the document, so we can do a direct\n// search-and-replace for each occurrence.\n\nconst replacements = [\n  [\"2025-11-09 Sunday\", \"2025-11-10 Monday\"],\n  [\"41\u00d788=\", \"31\u00d729=\"],\n  [\"44\u00d745=\", \"35\u00d768=\"],\n  [\"62\u00d765=\", \"77\u00d712=\"],\n  [\"31\u00d738=\", \"77\u00d733=\"],\n  [\"28\u00d787=\", \"25\u00d712=\"],\n  [\"57\u00d760=\", \"96\u00d752=\"],\n  [\"35\u00d712=\", \"15\u00d751=\"],\n  [\"76\u00d722=\", \"85\u00d753=\"],\n  [\"86\u00d712=\", \"34\u00d735=\"],\n  [\"91\u00d735=\", \"74\u00d729=\"],\n  [\"28\u00d752=\", \"70\u00d762=\"],\n  [\"27\u00d732=\", \"24\u00d780=\"],\n  [\"48\u00d775=\", \"54\u00d751=\"],\n  [\"80\u00d727=\", \"43\u00d785=\"],\n  [\"12\u00d768=\", \"19\u00d742=\"],\n  [\"79\u00d741=\", \"65\u00d789=\"],\n  [\"44\u00d784=\", \"65\u00d734=\"],\n  [\"36\u00d775=\", \"87\u00d749=\"],\n  [\"66\u00d756=\", \"46\u00d754=\"],\n  [\"80\u00d735=\", \"66\u00d782=\"],\n  [\"12\u00d775=\", \"53\u00d715=\"],\n  [\"27\u00d780=\", \"43\u00d715=\"],\n  [\"87\u00d780=\", \"44\u00d715=\"],\n  [\"12\u00d770=\", \"20\u00d766=\"],\n  [\"69\u00d757=\", \"27\u00d742=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date label and all multiplication problems in the table.\n# Each source text is unique in the document, so a simple Find/Replace\n# (scoped to the whole document) for each pair is sufficient.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-11-09 Sunday\", \"2025-11-10 Monday\"),\n    @(\"41\u00d788=\", \"31\u00d729=\"),\n    @(\"44\u00d745=\", \"35\u00d768=\"),\n    @(\"62\u00d765=\", \"77\u00d712=\"),\n    @(\"31\u00d738=\", \"77\u00d733=\"),\n    @(\"28\u00d787=\", \"25\u00d712=\"),\n    @(\"57\u00d760=\", \"96\u00d752=\"),\n    @(\"35\u00d712=\", \"15\u00d751=\"),\n    @(\"76\u00d722=\", \"85\u00d753=\"),\n    @(\"86\u00d712=\", \"34\u00d735=\"),\n    @(\"91\u00d735=\", \"74\u00d729=\"),\n    @(\"28\u00d752=\", \"70\u00d762=\"),\n    @(\"27\u00d732=\", \"24\u00d780=\"),\n    @(\"48\u00d775=\", \"54\u00d751=\"),\n    @(\"80\u00d727=\", \"43\u00d785=\"),\n    @(\"12\u00d768=\", \"19\u00d742=\"),\n    @(\"79\u00d741=\", \"65\u00d789=\"),\n    @(\"44\u00d784=\", \"65\u00d734=\"),\n    @(\"36\u00d775=\", \"87\u00d749=\"),\n    @(\"66\u00d756=\", \"46\u00d754=\"),\n    @(\"80\u00d735=\", \"66\u00d782=\"),\n    @(\"12\u00d775=\", \"53\u00d715=\"),\n    @(\"27\u00d780=\", \"43\u00d715=\"),\n    @(\"87\u00d780=\", \"44\u00d715=\"),\n    @(\"12\u00d770=\", \"20\u00d766=\"),\n    @(\"69\u00d757=\", \"27\u00d742=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
